$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the formatting of the last existing data row (279) down onto the
# three new rows so the new cells pick up the same styles (date format,
# border, bold header-like font on column A, etc.) without creating new
# style entries.
$ws.Range("A279:G279").Copy()
$ws.Range("A280:G282").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$data = @(
    @(44986.45833333334, "ECONOMICS:NGM2", 54191661320000, 54191661320000, 54191661320000, 54191661320000, 0),
    @(45017.45833333334, "ECONOMICS:NGM2", 55646745840000, 55646745840000, 55646745840000, 55646745840000, 0),
    @(45047.41666666666, "ECONOMICS:NGM2", 55500913410000, 55500913410000, 55500913410000, 55500913410000, 0)
)

$startRow = 280
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
